$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release History")

# New release row (row 19) for version 7.0.4
$ws.Range("A19").Value = "V3 EVT Firmware"
$ws.Range("B19").Value = "7.0.4"
$ws.Range("C19").Value = Get-Date -Year 2022 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0
$ws.Range("D19").Value = "Zound_Hendrix_M_Lite_V3_hwEVT_btswv7.0.4_20220415"
$ws.Range("E19").Value = "7.0.4"
$ws.Range("G19").Value = 0.6
$ws.Range("H19").Value = 3.1
$ws.Range("J19").Value = "N/A"
$ws.Range("K19").Value = "Resolve the recoonection issues."
